$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2: change Existencias (D2) value
$ws.Range("D2").Value = 4

# Row 3: becomes "Coca cola" / 1 / "Bebida" / 8
$ws.Range("A3").Value = "Coca cola"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Bebida"
$ws.Range("D3").Value = 8

# Row 4 (new): "Gorro de baño" / 2 / "Otros" / 9
$ws.Range("A4").Value = "Gorro de baño"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "Otros"
$ws.Range("D4").Value = 9

# Row 5 (new): "Chochos con tostado" / 1 / "Comida" / 15
$ws.Range("A5").Value = "Chochos con tostado"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "Comida"
$ws.Range("D5").Value = 15
